# "SF data through 2020-12-17"
#
# 1. Data sheet: one more day of observations (2020-12-17 / serial 44182)
#    becomes available and is filled into row 274 (the first previously-
#    blank row after 2020-12-16).
# 2. Internal sheet: a new Stan-sampler parameter ("warmup") is inserted
#    just above "cores", and "iter"/"adapt_delta" are retuned.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Data sheet
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Range("A274").Value = 44182
$wsData.Range("B274").Value = 161
$wsData.Range("C274").Value = 9
$wsData.Range("D274").Value = 41
$wsData.Range("E274").Value = 2
$wsData.Range("F274").Value = 172

$null = $wsData.Range("F275").Select()

# ---------------------------------------------------------------------
# Internal sheet
# ---------------------------------------------------------------------
$wsInt = $wb.Worksheets.Item("Internal")

# Insert a new row above row 9 ("cores") for the "warmup" parameter; the
# new row inherits the formatting of the row above it (row 8, "iter").
$wsInt.Rows.Item(9).Insert()
$wsInt.Range("A9").Value = "warmup"
$wsInt.Range("B9").Value = 1000

# iter: 2000 -> 1500
$wsInt.Range("B8").Value = 1500

# adapt_delta (now on row 13 after the insert above): 0.8 -> 0.9
$wsInt.Range("B13").Value = 0.9

$wsInt.Range("A10").Select()

# The "Internal" sheet becomes the active tab (was "Parameters with
# Distributions").
$wsInt.Activate()
